{"js": "// Update the \"Version 2.10.3\" line to read \"Version 2.13\" in bold, 10pt\n// text, remove the old \"Version 1.32\" line (and the blank lines that\n// used to surround it), and add four blank lines after the bookmark\n// paragraph near the end of the document (matching the target revision).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that currently reads \"Version 2.10.3\" and the one\n// that reads \"Version 1.32\" by scanning the paragraph text instead of\n// relying on fixed indices, so the script is resilient to minor layout\n// differences.\nlet versionParaIndex = -1;\nlet oldVersionParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Version 2.10.3\") !== -1) {\n    versionParaIndex = i;\n  } else if (text.indexOf(\"Version 1.32\") !== -1) {\n    oldVersionParaIndex = i;\n  }\n}\n\nif (versionParaIndex === -1) {\n  throw new Error(\"Could not find the 'Version 2.10.3' paragraph.\");\n}\n\n// Rewrite the version paragraph text/formatting: \"Version 2.13\", bold, 10pt.\nconst versionPara = paragraphs.items[versionParaIndex];\nversionPara.clear();\nversionPara.insertText(\"Version 2.13\", Word.InsertLocation.replace);\nversionPara.font.bold = true;\nversionPara.font.size = 10;\nversionPara.font.sizeBidirectional = 10;\nawait context.sync();\n\n// Remove the old \"Version 1.32\" paragraph plus the blank paragraphs\n// immediately before and after it (two blanks on each side).\nif (oldVersionParaIndex !== -1) {\n  const first = oldVersionParaIndex - 2;\n  const last = oldVersionParaIndex + 2;\n  for (let i = last; i >= first; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n\n// Insert four new blank paragraphs (matching the existing small/8pt\n// formatting) right after the paragraph that hosts the \"_GoBack\" bookmark.\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\n\nif (!bookmarkRange.isNullObject) {\n  let anchorPara = bookmarkRange.paragraphs.getFirst();\n  await context.sync();\n  for (let i = 0; i < 4; i++) {\n    anchorPara = anchorPara.insertParagraph(\"\", Word.InsertLocation.after);\n    anchorPara.font.size = 8;\n    anchorPara.font.sizeBidirectional = 8;\n    await context.sync();\n  }\n}\n", "ps1": "# Update the \"Version 2.10.3\" line to read \"Version 2.13\" in bold, 10pt\n# text, remove the old \"Version 1.32\" line (and the blank lines that used\n# to surround it), and add four blank lines after the bookmark paragraph\n# near the end of the document (matching the target revision).\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Locate the paragraph that currently reads \"Version 2.10.3\" and the one\n# that reads \"Version 1.32\" by scanning the paragraph text instead of\n# relying on fixed indices, so the script is resilient to minor layout\n# differences.\n$versionParaIndex = -1\n$oldVersionParaIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $t = $paras.Item($i).Range.Text\n  if ($t.Contains(\"Version 2.10.3\")) {\n    $versionParaIndex = $i\n  } elseif ($t.Contains(\"Version 1.32\")) {\n    $oldVersionParaIndex = $i\n  }\n}\n\n# Rewrite the version paragraph text/formatting: \"Version 2.13\", bold, 10pt.\nif ($versionParaIndex -ne -1) {\n  $p = $paras.Item($versionParaIndex)\n  $r = $p.Range\n  $r.End = $r.End - 1\n  $r.Text = \"Version 2.13\"\n\n  $full = $paras.Item($versionParaIndex).Range\n  $full.Font.Bold = 1\n  $full.Font.Size = 10\n  $full.Font.SizeBi = 10\n}\n\n# Remove the old \"Version 1.32\" paragraph plus the blank paragraphs\n# immediately before and after it (two blanks on each side).\nif ($oldVersionParaIndex -ne -1) {\n  $first = $oldVersionParaIndex - 2\n  $last = $oldVersionParaIndex + 2\n  for ($i = $last; $i -ge $first; $i--) {\n    $paras.Item($i).Range.Delete()\n  }\n}\n\n# Insert four new blank paragraphs (matching the existing small/8pt\n# formatting) right after the paragraph that hosts the \"_GoBack\" bookmark.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$anchor = $bm.Range.Paragraphs.Item(1)\nfor ($i = 0; $i -lt 4; $i++) {\n  $anchor.Range.InsertParagraphAfter()\n  $newRange = $anchor.Range.Next(4, 1).Range\n  $newRange.Font.Size = 8\n  $newRange.Font.SizeBi = 8\n  $anchor = $newRange.Paragraphs.Item(1)\n}\n"}
